$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to be treated as text so Excel doesn't silently coerce
    # numeric-looking strings (e.g. "599.95") into real numbers and lose
    # formatting (trailing zeros, etc.).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.291.24"
$ws.Range("E2").Value = "  -1.37%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.982.58"
$ws.Range("E3").Value = "  -0.19%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
Set-TextValue "D5" "599.95"
$ws.Range("E5").Value = "  +3.50%  "

# Row 6 - Solana
Set-TextValue "D6" "143.47"
$ws.Range("E6").Value = "  -1.53%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - XRP
Set-TextValue "D8" "0.518"
$ws.Range("E8").Value = "  -0.73%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.979.55"
$ws.Range("E9").Value = "  -0.33%  "

# Row 10 - Toncoin
Set-TextValue "D10" "6.05"
$ws.Range("E10").Value = "  +7.53%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -1.15%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000227"
$ws.Range("E13").Value = "  +0.27%  "

# Row 14 - Avalanche
Set-TextValue "D14" "34.25"
$ws.Range("E14").Value = "  -0.49%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +2.49%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "3.479.56"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  -2.02%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "61.312.96"
$ws.Range("E18").Value = "  -1.41%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "2.985.09"
$ws.Range("E19").Value = "  -0.26%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "449.03"
$ws.Range("E20").Value = "  -1.08%  "

# Row 21 - Chainlink
Set-TextValue "D21" "14.15"
$ws.Range("E21").Value = "  +2.32%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.683"
$ws.Range("E22").Value = "  +1.02%  "

# Row 23 - Uniswap
Set-TextValue "D23" "7.31"
$ws.Range("E23").Value = "  +0.63%  "

# Row 24 - Litecoin
Set-TextValue "D24" "81.90"
$ws.Range("E24").Value = "  +2.53%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  -2.93%  "

# Row 26 - RenderToken
Set-TextValue "D26" "10.48"
$ws.Range("E26").Value = "  +5.07%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "11.94"
$ws.Range("E27").Value = "  -2.45%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  +0.20%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +3.59%  "

# Row 30 - FirstDigitalUSD
$ws.Range("E30").Value = "  +0.02%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "7.14"
$ws.Range("E31").Value = "  +0.15%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  -1.42%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "27.15"
$ws.Range("E33").Value = "  +1.47%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.108"
$ws.Range("E34").Value = "  +1.84%  "

# Row 35 - PEPE
$ws.Range("D35").Value = "0.0₃0822"
$ws.Range("E35").Value = "  +5.17%  "

# Row 36 - Mantle
$ws.Range("E36").Value = "  +0.26%  "

# Row 37 - Filecoin
$ws.Range("E37").Value = "  +0.85%  "

# Row 38 - OKB
Set-TextValue "D38" "50.33"
$ws.Range("E38").Value = "  +0.41%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  -2.31%  "

# Row 40 - Cosmos
$ws.Range("E40").Value = "  +0.93%  "

# Row 41 & 42 - dogwifhat/Kaspa swap positions with updated values
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D41" "0.122"
$ws.Range("E41").Value = "  +10.12%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D42" "2.87"
$ws.Range("E42").Value = "  -1.73%  "

# Row 43 - Bittensor
Set-TextValue "D43" "397.24"
$ws.Range("E43").Value = "  -2.49%  "

# Row 44 - Arweave
Set-TextValue "D44" "39.70"
$ws.Range("E44").Value = "  +4.81%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  +0.34%  "

# Row 46 - TheGraph
$ws.Range("E46").Value = "  -2.10%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.688.18"
$ws.Range("E47").Value = "  -2.74%  "

# Row 48 - Monero
Set-TextValue "D48" "131.37"
$ws.Range("E48").Value = "  +2.59%  "

# Row 49 - USDe
$ws.Range("E49").Value = "  +0.10%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  -0.49%  "

# Row 51 - ThetaToken
$ws.Range("E51").Value = "  +0.91%  "
